# Insert a new weekly data row at position 323, shifting the existing
# rows 323:403 down to 324:404 (this also updates the sheet dimension
# automatically, same as Excel's native "Insert Row" behaviour).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("323:323").Insert()

# Populate the newly-inserted row 323 with the new observation.
$ws.Range("A323").Value2 = 6
$ws.Range("B323").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C323").Value2 = "Metropolitana"
$ws.Range("D323").Value2 = 44841
$ws.Range("E323").Value2 = 13
$ws.Range("F323").Value2 = "Fruta"
$ws.Range("G323").Value2 = 100101
$ws.Range("H323").Value2 = "Berries"
$ws.Range("I323").Value2 = 100101001
$ws.Range("J323").Value2 = "Arándano (blue)"
$ws.Range("K323").Value2 = "Sin especificar"
$ws.Range("L323").Value2 = "Primera"
$ws.Range("M323").Value2 = 2000
$ws.Range("N323").Value2 = 12000
$ws.Range("O323").Value2 = 12000
$ws.Range("P323").Value2 = 12000
$ws.Range("Q323").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R323").Value2 = "Región de O'Higgins"
$ws.Range("S323").Value2 = 6000
$ws.Range("T323").Value2 = 2
